# Add three new daily-snapshot rows (62-64) to each MID sheet,
# covering 2025-07-28 .. 2025-07-30, and extend the used range.
# (csv module error-handling fix upstream produced 3 extra rows of
# telemetry that were previously dropped.)
$wb = $excel.ActiveWorkbook

# Sheet 1: MID_LFT_#1
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(62, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 1).Value = 45848.46532407407
$ws.Cells.Item(62, 2).Value = "0x01,0x90"
$ws.Cells.Item(62, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(62, 4).Value = "0x01,0x5C"
$ws.Cells.Item(62, 5).Value = "0x07"
$ws.Cells.Item(62, 6).Value = 400
$ws.Cells.Item(62, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(62, 8).Value = 348
$ws.Cells.Item(62, 9).Value = 7

$ws.Cells.Item(63, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 1).Value = 45849.46129629629
$ws.Cells.Item(63, 2).Value = "0x01,0x90"
$ws.Cells.Item(63, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(63, 4).Value = "0x01,0x5C"
$ws.Cells.Item(63, 5).Value = "0x07"
$ws.Cells.Item(63, 6).Value = 400
$ws.Cells.Item(63, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(63, 8).Value = 348
$ws.Cells.Item(63, 9).Value = 7

$ws.Cells.Item(64, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 1).Value = 45850.46201388889
$ws.Cells.Item(64, 2).Value = "0x01,0x90"
$ws.Cells.Item(64, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1b,0x41,0x0c,"
$ws.Cells.Item(64, 4).Value = "0x01,0x58"
$ws.Cells.Item(64, 5).Value = "0x07"
$ws.Cells.Item(64, 6).Value = 400
$ws.Cells.Item(64, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(64, 8).Value = 344
$ws.Cells.Item(64, 9).Value = 7


# Sheet 2: MID_LFT_#2
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(62, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 1).Value = 45848.46532407407
$ws.Cells.Item(62, 2).Value = "0x01,0x7c"
$ws.Cells.Item(62, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(62, 4).Value = "0x01,0x58"
$ws.Cells.Item(62, 5).Value = "0x19"
$ws.Cells.Item(62, 6).Value = 380
$ws.Cells.Item(62, 7).Value = 568432987514711000000000.0
$ws.Cells.Item(62, 8).Value = 344
$ws.Cells.Item(62, 9).Value = 25

$ws.Cells.Item(63, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 1).Value = 45849.46129629629
$ws.Cells.Item(63, 2).Value = "0x01,0x7c"
$ws.Cells.Item(63, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(63, 4).Value = "0x01,0x58"
$ws.Cells.Item(63, 5).Value = "0x19"
$ws.Cells.Item(63, 6).Value = 380
$ws.Cells.Item(63, 7).Value = 568432987514711000000000.0
$ws.Cells.Item(63, 8).Value = 344
$ws.Cells.Item(63, 9).Value = 25

$ws.Cells.Item(64, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 1).Value = 45850.46201388889
$ws.Cells.Item(64, 2).Value = "0x01,0x7c"
$ws.Cells.Item(64, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c,"
$ws.Cells.Item(64, 4).Value = "0x01,0x54"
$ws.Cells.Item(64, 5).Value = "0x19"
$ws.Cells.Item(64, 6).Value = 380
$ws.Cells.Item(64, 7).Value = 568432987514711000000000.0
$ws.Cells.Item(64, 8).Value = 340
$ws.Cells.Item(64, 9).Value = 25


# Sheet 3: MID_PLT_#1
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(62, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 1).Value = 45848.46532407407
$ws.Cells.Item(62, 2).Value = "0x00,0x6e"
$ws.Cells.Item(62, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(62, 4).Value = "0x00,0x67"
$ws.Cells.Item(62, 5).Value = "0x15"
$ws.Cells.Item(62, 6).Value = 110
$ws.Cells.Item(62, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(62, 8).Value = 103
$ws.Cells.Item(62, 9).Value = 15

$ws.Cells.Item(63, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 1).Value = 45849.46129629629
$ws.Cells.Item(63, 2).Value = "0x00,0x6e"
$ws.Cells.Item(63, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(63, 4).Value = "0x00,0x67"
$ws.Cells.Item(63, 5).Value = "0x15"
$ws.Cells.Item(63, 6).Value = 110
$ws.Cells.Item(63, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(63, 8).Value = 103
$ws.Cells.Item(63, 9).Value = 15

$ws.Cells.Item(64, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 1).Value = 45850.46201388889
$ws.Cells.Item(64, 2).Value = "0x00,0x6e"
$ws.Cells.Item(64, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c,"
$ws.Cells.Item(64, 4).Value = "0x00,0x67"
$ws.Cells.Item(64, 5).Value = "0x15"
$ws.Cells.Item(64, 6).Value = 110
$ws.Cells.Item(64, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(64, 8).Value = 103
$ws.Cells.Item(64, 9).Value = 15


# Sheet 4: MID_PLT_#2
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(62, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 1).Value = 45848.46532407407
$ws.Cells.Item(62, 2).Value = "0x00,0x82"
$ws.Cells.Item(62, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(62, 4).Value = "0x00,0x7C"
$ws.Cells.Item(62, 5).Value = "0x9"
$ws.Cells.Item(62, 6).Value = 130
$ws.Cells.Item(62, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(62, 8).Value = 124
$ws.Cells.Item(62, 9).Value = 9

$ws.Cells.Item(63, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(63, 1).Value = 45849.46129629629
$ws.Cells.Item(63, 2).Value = "0x00,0x82"
$ws.Cells.Item(63, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(63, 4).Value = "0x00,0x7C"
$ws.Cells.Item(63, 5).Value = "0x9"
$ws.Cells.Item(63, 6).Value = 130
$ws.Cells.Item(63, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(63, 8).Value = 124
$ws.Cells.Item(63, 9).Value = 9

$ws.Cells.Item(64, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(64, 1).Value = 45850.46201388889
$ws.Cells.Item(64, 2).Value = "0x00,0x82"
$ws.Cells.Item(64, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c,"
$ws.Cells.Item(64, 4).Value = "0x00,0x7B"
$ws.Cells.Item(64, 5).Value = "0x9"
$ws.Cells.Item(64, 6).Value = 130
$ws.Cells.Item(64, 7).Value = 568631262647113000000000.0
$ws.Cells.Item(64, 8).Value = 123
$ws.Cells.Item(64, 9).Value = 9
